# Update on mail template
# - rename "AccountType" header to "accountType"
# - add a new "bankCode" column (H) with per-row bank codes
# - replace the row-2 sample record (bank / account number / account name)
#   with a new Stanbic IBTC Bank example

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dept Data")

# 1) Header rename: AccountType -> accountType
$ws.Range("G1").Value = "accountType"

# 2) New header for the added column
$ws.Range("H1").Value = "bankCode"

# 3) Row 2 gets a new sample record
$ws.Range("C2").Value = "Stanbic IBTC Bank"

# accountNumber needs to stay textual so the leading zero survives
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0034551560"

$ws.Range("E2").Value = "Prince Emmanuel Odu"

# 4) Populate the new bankCode column for every existing data row
$ws.Range("H2").Value = 221
$ws.Range("H3").Value = 21
$ws.Range("H4").Value = 200
$ws.Range("H5").Value = 341
$ws.Range("H6").Value = 245
$ws.Range("H7").Value = 200
